$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).ColumnWidth = 12.109375
$ws.Columns.Item(3).ColumnWidth = 17.5546875
$ws.Columns.Item(6).ColumnWidth = 22.6640625

$ws.Range("B3:B6").HorizontalAlignment = -4108
$ws.Range("F3:F6").HorizontalAlignment = -4108

$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0

$ws.Range("B7").Font.Bold = $true
$ws.Range("B7").Interior.ThemeColor = 7
$ws.Range("B7").Interior.TintAndShade = 0.59999389629810485
$ws.Range("B7").HorizontalAlignment = -4108

$ws.Range("C7:E7").Interior.ThemeColor = 7
$ws.Range("C7:E7").Interior.TintAndShade = 0.59999389629810485
$ws.Range("D7").Value = -50
$ws.Range("E7").Value = 30

$ws.Range("F3:F6").Select()

$wb.Worksheets(1).ListObjects(1).ListColumns(1).Range.Select()
